$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date cells (I1/K1 text dates, M1 numeric "Run Date" serial)
$ws.Range("I1").Value2 = "06-11-2025 00:00:00"
$ws.Range("K1").Value2 = "06-11-2025 00:00:00"
$ws.Range("M1").Value2 = 45967

# Cyclically rotate B/D/E/F/G values among each block of rows (each row
# takes on the next row's original values, wrapping around).
$blocks = @(
    @(136,137),
    @(146,147,148),
    @(163,164),
    @(233,234),
    @(246,247),
    @(277,278),
    @(292,293),
    @(295,296),
    @(311,312),
    @(420,421),
    @(465,466),
    @(472,473),
    @(490,491),
    @(596,597),
    @(705,706),
    @(732,733)
)

foreach ($block in $blocks) {
    $n = $block.Length
    $origB = @()
    $origD = @()
    $origE = @()
    $origF = @()
    $origG = @()
    foreach ($r in $block) {
        $origB += $ws.Cells.Item($r, 2).Value2
        $origD += $ws.Cells.Item($r, 4).Value2
        $origE += $ws.Cells.Item($r, 5).Value2
        $origF += $ws.Cells.Item($r, 6).Value2
        $origG += $ws.Cells.Item($r, 7).Value2
    }
    for ($i = 0; $i -lt $n; $i++) {
        $next = ($i + 1) % $n
        $r = $block[$i]
        $ws.Cells.Item($r, 2).Value2 = $origB[$next]
        $ws.Cells.Item($r, 4).Value2 = $origD[$next]
        $ws.Cells.Item($r, 5).Value2 = $origE[$next]
        $ws.Cells.Item($r, 6).Value2 = $origF[$next]
        $ws.Cells.Item($r, 7).Value2 = $origG[$next]
    }
}
